$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-22 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-23 Friday", 2) | Out-Null
$d.Content.Find.Execute("464÷8=58, 0", $true, $false, $false, $false, $false, $true, 1, $false, "522÷8=65, 2", 2) | Out-Null
$d.Content.Find.Execute("708÷4=177, 0", $true, $false, $false, $false, $false, $true, 1, $false, "862÷5=172, 2", 2) | Out-Null
$d.Content.Find.Execute("552÷7=78, 6", $true, $false, $false, $false, $false, $true, 1, $false, "820÷3=273, 1", 2) | Out-Null
$d.Content.Find.Execute("782÷9=86, 8", $true, $false, $false, $false, $false, $true, 1, $false, "497÷9=55, 2", 2) | Out-Null
$d.Content.Find.Execute("453÷6=75, 3", $true, $false, $false, $false, $false, $true, 1, $false, "866÷7=123, 5", 2) | Out-Null
$d.Content.Find.Execute("182÷5=36, 2", $true, $false, $false, $false, $false, $true, 1, $false, "196÷6=32, 4", 2) | Out-Null
$d.Content.Find.Execute("264÷5=52, 4", $true, $false, $false, $false, $false, $true, 1, $false, "484÷4=121, 0", 2) | Out-Null
$d.Content.Find.Execute("677÷5=135, 2", $true, $false, $false, $false, $false, $true, 1, $false, "340÷9=37, 7", 2) | Out-Null
$d.Content.Find.Execute("574÷3=191, 1", $true, $false, $false, $false, $false, $true, 1, $false, "592÷2=296, 0", 2) | Out-Null
$d.Content.Find.Execute("190÷2=95, 0", $true, $false, $false, $false, $false, $true, 1, $false, "660÷6=110, 0", 2) | Out-Null
$d.Content.Find.Execute("466÷7=66, 4", $true, $false, $false, $false, $false, $true, 1, $false, "699÷5=139, 4", 2) | Out-Null
$d.Content.Find.Execute("628÷4=157, 0", $true, $false, $false, $false, $false, $true, 1, $false, "897÷8=112, 1", 2) | Out-Null
$d.Content.Find.Execute("544÷3=181, 1", $true, $false, $false, $false, $false, $true, 1, $false, "965÷9=107, 2", 2) | Out-Null
$d.Content.Find.Execute("842÷3=280, 2", $true, $false, $false, $false, $false, $true, 1, $false, "248÷6=41, 2", 2) | Out-Null
$d.Content.Find.Execute("762÷3=254, 0", $true, $false, $false, $false, $false, $true, 1, $false, "170÷6=28, 2", 2) | Out-Null
$d.Content.Find.Execute("284÷9=31, 5", $true, $false, $false, $false, $false, $true, 1, $false, "507÷4=126, 3", 2) | Out-Null
$d.Content.Find.Execute("336÷6=56, 0", $true, $false, $false, $false, $false, $true, 1, $false, "427÷4=106, 3", 2) | Out-Null
$d.Content.Find.Execute("159÷6=26, 3", $true, $false, $false, $false, $false, $true, 1, $false, "294÷6=49, 0", 2) | Out-Null
$d.Content.Find.Execute("845÷3=281, 2", $true, $false, $false, $false, $false, $true, 1, $false, "281÷2=140, 1", 2) | Out-Null
$d.Content.Find.Execute("127÷7=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "286÷4=71, 2", 2) | Out-Null
$d.Content.Find.Execute("149÷7=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "365÷2=182, 1", 2) | Out-Null
$d.Content.Find.Execute("535÷4=133, 3", $true, $false, $false, $false, $false, $true, 1, $false, "661÷3=220, 1", 2) | Out-Null
$d.Content.Find.Execute("840÷8=105, 0", $true, $false, $false, $false, $false, $true, 1, $false, "120÷9=13, 3", 2) | Out-Null
$d.Content.Find.Execute("875÷8=109, 3", $true, $false, $false, $false, $false, $true, 1, $false, "245÷5=49, 0", 2) | Out-Null
$d.Content.Find.Execute("158÷7=22, 4", $true, $false, $false, $false, $false, $true, 1, $false, "964÷2=482, 0", 2) | Out-Null
